# edit.ps1 - apply the "Fix funcionalidades Actualizacion caso 001 al 012 operativos" change
# to the DataPruebas sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "DataPruebas"

# ---------------------------------------------------------------------------
# Row 7 : swap the placeholder passenger/payment data for the new ones
# ---------------------------------------------------------------------------
$ws.Range("I7").Value = " "
$ws.Range("J7").Value = " "

# ---------------------------------------------------------------------------
# Row 10 : ticket price becomes a real currency number instead of plain text
# (applied before the V7 underline tweak below so the new style entries come
# out in the same order as the authoritative edit)
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 21.66
$ws.Range("B10").NumberFormat = "_ * #,##0.00_-\ [$€-1]_ ;_ * #,##0.00\-\ [$€-1]_ ;_ * ""-""??_-\ [$€-1]_ ;_ @_ "

# ---------------------------------------------------------------------------
# Row 11 : B11 gets a single-space placeholder value (was empty)
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = " "

# V7 used to hold a raw long card number; it is now a formatted card-number
# string and gets an underline applied to it.
$ws.Range("V7").Value = "5100 0100 0000 0114"
$ws.Range("V7").Font.Underline = $true

$ws.Range("R7").Value = "Messi"
$ws.Range("Q7").Value = "Lia"

# ---------------------------------------------------------------------------
# Rows 12-13 : same card-number text swap as row 7 (keeps existing style)
# ---------------------------------------------------------------------------
$ws.Range("L12").Value = "5100 0100 0000 0114"
$ws.Range("L13").Value = "5100 0100 0000 0114"

# Row 13 : small numeric correction
$ws.Range("N13").Value = 23

# ---------------------------------------------------------------------------
# Column widths: column A grew wider, and the new column H got an explicit
# width as well.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 54.5
$ws.Columns.Item(8).ColumnWidth = 9.92

# ---------------------------------------------------------------------------
# View: scroll position / selection moved from S17 to B7
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B7").Select()
